$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to be stored as text, matching the source data
# (inline strings), so purely numeric-looking values like "0.997" are not
# auto-converted into numeric cells by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "63.289.10"
$ws.Range("E2").Value = "  +6.56%  "

# Row 3
$ws.Range("D3").Value = "2.421.01"
$ws.Range("E3").Value = "  +2.51%  "

# Row 4
$ws.Range("E4").Value = "  +0.57%  "

# Row 5
$ws.Range("D5").Value = "576.27"
$ws.Range("E5").Value = "  +3.34%  "

# Row 6
$ws.Range("D6").Value = "146.28"
$ws.Range("E6").Value = "  +6.61%  "

# Row 7
$ws.Range("D7").Value = "0.997"
$ws.Range("E7").Value = "  -0.47%  "

# Row 8
$ws.Range("D8").Value = "0.541"
$ws.Range("E8").Value = "  +2.71%  "

# Row 9
$ws.Range("D9").Value = "2.453.50"
$ws.Range("E9").Value = "  +4.12%  "

# Row 10
$ws.Range("E10").Value = "  +6.35%  "

# Row 11
$ws.Range("E11").Value = "  +1.50%  "

# Row 12
$ws.Range("D12").Value = "5.25"
$ws.Range("E12").Value = "  +3.93%  "

# Row 13
$ws.Range("D13").Value = "0.354"
$ws.Range("E13").Value = "  +5.71%  "

# Row 14
$ws.Range("D14").Value = "27.13"
$ws.Range("E14").Value = "  +6.86%  "

# Row 15
$ws.Range("D15").Value = "0.0000179"
$ws.Range("E15").Value = "  +9.46%  "

# Row 16
$ws.Range("D16").Value = "2.862.49"
$ws.Range("E16").Value = "  +2.91%  "

# Row 17
$ws.Range("D17").Value = "62.987.31"
$ws.Range("E17").Value = "  +6.28%  "

# Row 18
$ws.Range("D18").Value = "2.459.41"
$ws.Range("E18").Value = "  +4.39%  "

# Row 19
$ws.Range("D19").Value = "7.96"
$ws.Range("E19").Value = "  +0.42%  "

# Row 20
$ws.Range("D20").Value = "10.98"
$ws.Range("E20").Value = "  +5.62%  "

# Row 21
$ws.Range("D21").Value = "328.66"
$ws.Range("E21").Value = "  +2.59%  "

# Row 22
$ws.Range("E22").Value = "  +2.91%  "

# Row 23
$ws.Range("E23").Value = "  +13.97%  "

# Row 24
$ws.Range("D24").Value = "0.999"
$ws.Range("E24").Value = "  -0.29%  "

# Row 25
$ws.Range("D25").Value = "66.02"
$ws.Range("E25").Value = "  +3.17%  "

# Row 26
$ws.Range("D26").Value = "621.00"
$ws.Range("E26").Value = "  +13.11%  "

# Row 27
$ws.Range("B27").Value = "Binance-PegBSC-USD"
$ws.Range("C27").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D27").Value = "1.10"
$ws.Range("E27").Value = "  +10.29%  "

# Row 28
$ws.Range("B28").Value = "Aptos"
$ws.Range("C28").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D28").Value = "8.44"
$ws.Range("E28").Value = "  +4.32%  "

# Row 29
$ws.Range("B29").Value = "PEPE"
$ws.Range("C29").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D29").Value = "0.0₃0988"
$ws.Range("E29").Value = "  +8.57%  "

# Row 30
$ws.Range("B30").Value = "WrappedeETH"
$ws.Range("C30").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D30").Value = "2.565.35"
$ws.Range("E30").Value = "  +3.82%  "

# Row 31
$ws.Range("B31").Value = "InternetComputer(DFINITY)"
$ws.Range("C31").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D31").Value = "8.19"
$ws.Range("E31").Value = "  +3.49%  "

# Row 32
$ws.Range("B32").Value = "Fetch.AI"
$ws.Range("C32").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D32").Value = "1.42"
$ws.Range("E32").Value = "  +9.98%  "

# Row 33
$ws.Range("D33").Value = "0.138"
$ws.Range("E33").Value = "  +6.68%  "

# Row 34
$ws.Range("B34").Value = "PancakeSwap"
$ws.Range("C34").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D34").Value = "1.85"
$ws.Range("E34").Value = "  +4.83%  "

# Row 35
$ws.Range("B35").Value = "ImmutableX"
$ws.Range("C35").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D35").Value = "1.49"
$ws.Range("E35").Value = "  +6.32%  "

# Row 36
$ws.Range("D36").Value = "0.995"
$ws.Range("E36").Value = "  -0.48%  "

# Row 37
$ws.Range("D37").Value = "4.76"
$ws.Range("E37").Value = "  +5.56%  "

# Row 38
$ws.Range("D38").Value = "0.374"
$ws.Range("E38").Value = "  +2.81%  "

# Row 39
$ws.Range("D39").Value = "152.80"
$ws.Range("E39").Value = "  +2.13%  "

# Row 40
$ws.Range("D40").Value = "5.43"
$ws.Range("E40").Value = "  +8.99%  "

# Row 41
$ws.Range("D41").Value = "18.69"
$ws.Range("E41").Value = "  +3.94%  "

# Row 42
$ws.Range("D42").Value = "2.76"
$ws.Range("E42").Value = "  +15.10%  "

# Row 43
$ws.Range("D43").Value = "1.77"
$ws.Range("E43").Value = "  +8.15%  "

# Row 44
$ws.Range("D44").Value = "42.38"
$ws.Range("E44").Value = "  +2.34%  "

# Row 45
$ws.Range("E45").Value = "  -0.03%  "

# Row 46
$ws.Range("B46").Value = "BabyDogeCoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D46").Value = "0.0₆0286"
$ws.Range("E46").Value = "  -0.43%  "

# Row 47
$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D47").Value = "144.87"
$ws.Range("E47").Value = "  +5.17%  "

# Row 48
$ws.Range("B48").Value = "Filecoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D48").Value = "3.60"
$ws.Range("E48").Value = "  +3.32%  "

# Row 49
$ws.Range("B49").Value = "InjectiveProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D49").Value = "20.34"
$ws.Range("E49").Value = "  +7.19%  "

# Row 50
$ws.Range("B50").Value = "Mantle"
$ws.Range("C50").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D50").Value = "0.603"
$ws.Range("E50").Value = "  +3.58%  "

# Row 51
$ws.Range("B51").Value = "Hedera"
$ws.Range("C51").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D51").Value = "0.0518"
$ws.Range("E51").Value = "  +4.22%  "
